$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 380, shifting rows 380:475 down to 382:477.
$ws.Rows("380:381").Insert()

# New row 380: Packham's Triumph / Primera, week of 2022-07-12
$ws.Cells.Item(380, 1).Value = 11
$ws.Cells.Item(380, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(380, 3).Value = "Bíobío"
$ws.Cells.Item(380, 4).Value = 44754
$ws.Cells.Item(380, 5).Value = 8
$ws.Cells.Item(380, 6).Value = "Fruta"
$ws.Cells.Item(380, 7).Value = 100104
$ws.Cells.Item(380, 8).Value = "Frutos de pepita"
$ws.Cells.Item(380, 9).Value = 100104005
$ws.Cells.Item(380, 10).Value = "Pera"
$ws.Cells.Item(380, 11).Value = "Packham's Triumph"
$ws.Cells.Item(380, 12).Value = "Primera"
$ws.Cells.Item(380, 13).Value = 250
$ws.Cells.Item(380, 14).Value = 8000
$ws.Cells.Item(380, 15).Value = 8500
$ws.Cells.Item(380, 16).Value = 8300
$ws.Cells.Item(380, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(380, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(380, 19).Value = 519
$ws.Cells.Item(380, 20).Value = 16

# New row 381: Packham's Triumph / Segunda, week of 2022-07-12
$ws.Cells.Item(381, 1).Value = 11
$ws.Cells.Item(381, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(381, 3).Value = "Bíobío"
$ws.Cells.Item(381, 4).Value = 44754
$ws.Cells.Item(381, 5).Value = 8
$ws.Cells.Item(381, 6).Value = "Fruta"
$ws.Cells.Item(381, 7).Value = 100104
$ws.Cells.Item(381, 8).Value = "Frutos de pepita"
$ws.Cells.Item(381, 9).Value = 100104005
$ws.Cells.Item(381, 10).Value = "Pera"
$ws.Cells.Item(381, 11).Value = "Packham's Triumph"
$ws.Cells.Item(381, 12).Value = "Segunda"
$ws.Cells.Item(381, 13).Value = 250
$ws.Cells.Item(381, 14).Value = 6500
$ws.Cells.Item(381, 15).Value = 7000
$ws.Cells.Item(381, 16).Value = 6740
$ws.Cells.Item(381, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(381, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(381, 19).Value = 421
$ws.Cells.Item(381, 20).Value = 16
